$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the "longtitude" -> "longitude" typo in the header row (B1)
$ws.Range("B1").Value = "longitude"

# Add the new vaccination center row (row 66), copying the style/number
# format from the row directly above it (row 65) so the new cells match
# the existing latitude/longitude/Barangay/Name formatting.
$ws.Range("A65:D65").Copy() | Out-Null
$ws.Range("A66:D66").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("A66").Value = 13.8280654994402
$ws.Range("B66").Value = 121.394678081324
$ws.Range("C66").Value = "Poblacion"
$ws.Range("D66").Value = "San Juan Nepomuceno Church"

# Match the refreshed view state: zoomed to 130% with D13 selected.
$excel.ActiveWindow.Zoom = 130
$ws.Range("D13").Select() | Out-Null
